$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-09 Wednesday" "2025-07-10 Thursday"

Replace-Text "44÷5=8, 4" "93÷7=13, 2"
Replace-Text "28÷4=7, 0" "77÷9=8, 5"
Replace-Text "20÷3=6, 2" "48÷7=6, 6"
Replace-Text "75÷5=15, 0" "29÷8=3, 5"
Replace-Text "21÷4=5, 1" "70÷2=35, 0"

Replace-Text "22÷7=3, 1" "25÷8=3, 1"
Replace-Text "27÷4=6, 3" "50÷6=8, 2"
Replace-Text "84÷4=21, 0" "33÷7=4, 5"
Replace-Text "31÷5=6, 1" "64÷4=16, 0"
Replace-Text "79÷4=19, 3" "93÷6=15, 3"

Replace-Text "64÷2=32, 0" "75÷9=8, 3"
Replace-Text "14÷3=4, 2" "10÷7=1, 3"
Replace-Text "47÷2=23, 1" "11÷4=2, 3"
Replace-Text "36÷2=18, 0" "42÷8=5, 2"
Replace-Text "59÷3=19, 2" "14÷8=1, 6"

Replace-Text "36÷7=5, 1" "32÷4=8, 0"
Replace-Text "40÷5=8, 0" "84÷3=28, 0"
Replace-Text "14÷6=2, 2" "91÷7=13, 0"
Replace-Text "21÷8=2, 5" "17÷2=8, 1"
Replace-Text "92÷2=46, 0" "77÷9=8, 5"

Replace-Text "92÷9=10, 2" "53÷5=10, 3"
Replace-Text "50÷7=7, 1" "69÷9=7, 6"
Replace-Text "88÷3=29, 1" "28÷3=9, 1"
Replace-Text "37÷2=18, 1" "82÷7=11, 5"
Replace-Text "24÷9=2, 6" "84÷4=21, 0"
